$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain text, matching source data semantics
# (values like "63.292.07" or "0.999" must not be auto-converted to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.292.07'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '3.225.47'
$ws.Range('E3').Value = '  +2.81%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '594.09'
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').Value = '140.94'
$ws.Range('E6').Value = '  -1.42%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '3.225.22'
$ws.Range('E8').Value = '  +2.60%  '
$ws.Range('D9').Value = '0.521'
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('D10').Value = '0.148'
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('D11').Value = '5.39'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '0.465'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  -2.06%  '
$ws.Range('D14').Value = '34.49'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = '3.715.79'
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '3.193.95'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '63.275.90'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('D19').Value = '6.76'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').Value = '478.56'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').Value = '14.13'
$ws.Range('E21').Value = '  -3.74%  '
$ws.Range('D22').Value = '0.713'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').Value = '7.86'
$ws.Range('E23').Value = '  +2.95%  '
$ws.Range('D24').Value = '84.28'
$ws.Range('E24').Value = '  -3.21%  '
$ws.Range('D25').Value = '13.19'
$ws.Range('E25').Value = '  -1.78%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '2.73'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').Value = '7.23'
$ws.Range('E28').Value = '  +3.37%  '
$ws.Range('D29').Value = '8.07'
$ws.Range('E29').Value = '  -2.06%  '
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  +2.82%  '
$ws.Range('D31').Value = '27.43'
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('D33').Value = '0.109'
$ws.Range('E33').Value = '  -1.87%  '
$ws.Range('D34').Value = '2.56'
$ws.Range('E34').Value = '  -2.85%  '
$ws.Range('D35').Value = '1.09'
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('D36').Value = '5.87'
$ws.Range('E36').Value = '  -2.40%  '
$ws.Range('D37').Value = '52.79'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').Value = '0.0₃0720'
$ws.Range('E38').Value = '  -3.20%  '
$ws.Range('D39').Value = '0.0393'
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('D40').Value = '424.33'
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.79'
$ws.Range('E41').Value = '  -6.15%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Value = '8.41'
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('D43').Value = '2.980.36'
$ws.Range('E43').Value = '  +3.42%  '
$ws.Range('D44').Value = '0.112'
$ws.Range('E44').Value = '  -7.24%  '
$ws.Range('D45').Value = '0.268'
$ws.Range('E45').Value = '  +2.88%  '
$ws.Range('D46').Value = '2.16'
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '25.91'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').Value = '2.28'
$ws.Range('E50').Value = '  -4.60%  '
$ws.Range('D51').Value = '119.95'
$ws.Range('E51').Value = '  -0.78%  '
